$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the touched cells so values like "305.36", "0.00%"
# round-trip as exact text (matching the source inlineStr cells) instead of
# being auto-coerced into numbers/percentages by Excel.
$cells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'B7', 'C7', 'D7', 'E7', 'B8', 'C8', 'D8', 'E8', 'B9', 'C9', 'D9', 'E9', 'B10', 'C10', 'D10', 'E10', 'B11', 'C11', 'D11', 'E11', 'B12', 'C12', 'D12', 'E12', 'B13', 'C13', 'D13', 'E13', 'B14', 'C14', 'D14', 'E14', 'B15', 'C15', 'D15', 'E15', 'B16', 'C16', 'D16', 'E16', 'B17', 'C17', 'D17', 'E17', 'E18', 'D19', 'E19', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'D26', 'E26', 'E27', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'E47', 'E48', 'D49', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '305.36'
$ws.Range('E2').Value = '0.00%'
$ws.Range('D3').Value = '36.33'
$ws.Range('E3').Value = '-1.47%'
$ws.Range('D4').Value = '5.059'
$ws.Range('E4').Value = '0.66%'
$ws.Range('D5').Value = '0.07847'
$ws.Range('E5').Value = '0.10%'
$ws.Range('D6').Value = '2.316'
$ws.Range('E6').Value = '6.96%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '4.158'
$ws.Range('E7').Value = '2.47%'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').Value = '7.986'
$ws.Range('E8').Value = '-0.63%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '0.9295'
$ws.Range('E9').Value = '0.65%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '0.09752'
$ws.Range('E10').Value = '-1.94%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1865'
$ws.Range('E11').Value = '-0.39%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.08903'
$ws.Range('E12').Value = '2.19%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.03789'
$ws.Range('E13').Value = '5.65%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.09905'
$ws.Range('E14').Value = '-0.36%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '0.001435'
$ws.Range('E15').Value = '-3.02%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.005678'
$ws.Range('E16').Value = '0.15%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.469'
$ws.Range('E17').Value = '0.09%'
$ws.Range('E18').Value = '13.98%'
$ws.Range('D19').Value = '0.3421'
$ws.Range('E19').Value = '-0.77%'
$ws.Range('E20').Value = '-1.65%'
$ws.Range('D21').Value = '5.138'
$ws.Range('E21').Value = '4.59%'
$ws.Range('D22').Value = '0.2261'
$ws.Range('E22').Value = '2.74%'
$ws.Range('D23').Value = '0.04599'
$ws.Range('E23').Value = '-0.04%'
$ws.Range('D24').Value = '0.001237'
$ws.Range('E24').Value = '0.41%'
$ws.Range('D25').Value = '0.004760'
$ws.Range('D26').Value = '0.0001308'
$ws.Range('E26').Value = '-6.62%'
$ws.Range('E27').Value = '74.30%'
$ws.Range('D39').Value = '0.01928'
$ws.Range('E39').Value = '6.29%'
$ws.Range('D40').Value = '0.05166'
$ws.Range('E40').Value = '8.89%'
$ws.Range('D41').Value = '0.007812'
$ws.Range('E41').Value = '-1.01%'
$ws.Range('D42').Value = '0.1386'
$ws.Range('E42').Value = '-1.40%'
$ws.Range('D43').Value = '0.007849'
$ws.Range('E43').Value = '3.35%'
$ws.Range('D44').Value = '0.002153'
$ws.Range('E44').Value = '-3.94%'
$ws.Range('D45').Value = '0.01126'
$ws.Range('E45').Value = '7.41%'
$ws.Range('D46').Value = '0.00006183'
$ws.Range('E46').Value = '-2.37%'
$ws.Range('E47').Value = '0.55%'
$ws.Range('E48').Value = '0.14%'
$ws.Range('D49').Value = '51.70'
$ws.Range('E49').Value = '54.34%'
$ws.Range('D50').Value = '0.001911'
$ws.Range('E50').Value = '-29.03%'
$ws.Range('D51').Value = '0.00002113'
$ws.Range('E51').Value = '0.55%'
